$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.475135326385498
$ws.Range("B1").Value = 2.97974681854248
$ws.Range("C1").Value = 2.631460666656494
$ws.Range("D1").Value = 2.395705223083496
$ws.Range("E1").Value = 1.720226645469666
